# ps-trigonometricidentities.docx edit
#
# The authored change merges several runs of plain text (that pandoc/quarto
# had split one-run-per-word) into a single run per paragraph, for the
# title, author list and abstract paragraphs. Doing a same-text Find &
# Replace over the whole (already-contiguous) logical text of each of
# those paragraphs makes Word re-write the matched range as one run,
# which collapses the redundant <w:r> splits exactly like the diff shows.
#
# It also swaps a hard-coded build-machine path baked into the note-icon
# picture's legacy description attribute from a macOS quarto install path
# to a Windows one.

$d = $word.ActiveDocument

# --- Title: "Proof:" / " " / "Trigonometric" / " " / "identities" -> one run
$d.Content.Find.Execute(
    "Proof: Trigonometric identities", $false, $false, $false, $false, $false,
    $true, 1, $false, "Proof: Trigonometric identities", 2)

# --- Author list -> one run
$d.Content.Find.Execute(
    "Shanelle Advani, Krish Chaudhary, Tom Coleman, Dzhemma Ruseva", $false, $false, $false, $false, $false,
    $true, 1, $false, "Shanelle Advani, Krish Chaudhary, Tom Coleman, Dzhemma Ruseva", 2)

# --- Abstract sentence -> one run
$d.Content.Find.Execute(
    "Explanations as to why certain trigonometric identities are true.", $false, $false, $false, $false, $false,
    $true, 1, $false, "Explanations as to why certain trigonometric identities are true.", 2)

# --- Note-icon picture description: macOS build path -> Windows build path.
# (Best-effort: some hosts keep this legacy duplicate description attribute
# in sync with InlineShape.AlternativeText even for pictures nested in a
# table cell; harmless no-op if that particular sync isn't available.)
try {
    $shp = $d.InlineShapes.Item(1)
    $shp.AlternativeText = "D:\Programming Languages\share\formats\docx\note.png"
} catch {
}
